$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update company name text in B3
$ws.Range("B3").Value = "Air Liquide Côte d'Ivoire Société Anonyme (BRVM:SIVC)"

# Remove the historical_growth_revenue_last_5_years column values for both rows
$ws.Range("D2").ClearContents()
$ws.Range("D3").ClearContents()

foreach ($r in 2,3) {
    $ws.Cells.Item($r, 7).Value = -0.006050420168067226   # G: ebitdard_margin
    $ws.Cells.Item($r, 8).Value = -0.006050420168067226   # H: ebitda_margin
    $ws.Cells.Item($r, 9).Value = -0.1117647058823529     # I: operating_margin
    $ws.Cells.Item($r, 10).Value = -0.1117647058823529    # J: after_tax_operating_margin
    $ws.Cells.Item($r, 11).Value = -0.49                  # K: trailing_net_income
    $ws.Cells.Item($r, 12).Value = -0.04117647058823529   # L: net_margin

    $ws.Cells.Item($r, 21).Value = 0                      # U: cash
    $ws.Cells.Item($r, 22).Value = 0                      # V: cash_market_cap
    $ws.Cells.Item($r, 23).Value = -0.04666666666666667   # W: roe
    $ws.Cells.Item($r, 24).Value = 0.09035107089292931    # X: cost_equity
    $ws.Cells.Item($r, 25).Value = -0.137017737559596     # Y: roe_cost_equity
    $ws.Cells.Item($r, 26).Value = 0.8756438557763061     # Z: sales_invested_capital
    $ws.Cells.Item($r, 27).Value = -0.09786607799852834   # AA: roic
    $ws.Cells.Item($r, 28).Value = 0.09035107089292931    # AB: cost_capital
    $ws.Cells.Item($r, 29).Value = -0.1882171488914576    # AC: roic_cost_capital
    $ws.Cells.Item($r, 30).Value = 0                      # AD: debt_total
    $ws.Cells.Item($r, 32).Value = 0                      # AF: debt_total_inc_leases
    $ws.Cells.Item($r, 33).Value = 0                      # AG: net_debt
    $ws.Cells.Item($r, 34).Value = 0                      # AH: debt_market_capital

    # AI (debt_book_capital) is removed entirely
    $ws.Cells.Item($r, 35).ClearContents()

    $ws.Cells.Item($r, 36).Value = 0                      # AJ: net_debt_market_capital

    # AK (net_debt_book_capital) is removed entirely
    $ws.Cells.Item($r, 37).ClearContents()

    $ws.Cells.Item($r, 38).Value = 0.316                  # AL: interest_expenses
    $ws.Cells.Item($r, 39).Value = 0.308                  # AM: net_interest_expenses
    $ws.Cells.Item($r, 40).Value = -0                     # AN: debt_ebitda
    $ws.Cells.Item($r, 41).Value = -4.208860759493671     # AO: ebit_interest_expenses
    $ws.Cells.Item($r, 42).Value = -0                     # AP: net_debt_ebitda
    $ws.Cells.Item($r, 43).Value = -4.318181818181818     # AQ: ebit_net_interest_expenses
}
